$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'42.813.15"
$ws.Cells.Item(2, 5).Value = "  +4.79%  "
$ws.Cells.Item(3, 4).Value = "'2.249.22"
$ws.Cells.Item(3, 5).Value = "  +3.84%  "
$ws.Cells.Item(4, 5).Value = "  +0.06%  "
$ws.Cells.Item(5, 4).Value = "'249.67"
$ws.Cells.Item(5, 5).Value = "  -0.03%  "
$ws.Cells.Item(6, 4).Value = "'0.633"
$ws.Cells.Item(6, 5).Value = "  +2.61%  "
$ws.Cells.Item(7, 4).Value = "'69.92"
$ws.Cells.Item(7, 5).Value = "  +4.86%  "
$ws.Cells.Item(8, 2).Value = "USDC"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Cells.Item(8, 4).Value = "'1.00"
$ws.Cells.Item(8, 5).Value = "  -0.12%  "
$ws.Cells.Item(9, 2).Value = "Cardano"
$ws.Cells.Item(9, 3).Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Cells.Item(9, 4).Value = "'0.669"
$ws.Cells.Item(9, 5).Value = "  +18.97%  "
$ws.Cells.Item(10, 4).Value = "'39.48"
$ws.Cells.Item(10, 5).Value = "  +11.08%  "
$ws.Cells.Item(11, 2).Value = "OKB"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(11, 4).Value = "'59.16"
$ws.Cells.Item(11, 5).Value = "  +1.65%  "
$ws.Cells.Item(12, 2).Value = "Dogecoin"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Cells.Item(12, 4).Value = "'0.0960"
$ws.Cells.Item(12, 5).Value = "  +4.38%  "
$ws.Cells.Item(13, 4).Value = "'7.58"
$ws.Cells.Item(13, 5).Value = "  +10.42%  "
$ws.Cells.Item(14, 5).Value = "  +1.42%  "
$ws.Cells.Item(15, 4).Value = "'2.586.35"
$ws.Cells.Item(15, 5).Value = "  +3.63%  "
$ws.Cells.Item(16, 2).Value = "Chainlink"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Cells.Item(16, 4).Value = "'14.77"
$ws.Cells.Item(16, 5).Value = "  +4.75%  "
$ws.Cells.Item(17, 2).Value = "Polygon"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(17, 4).Value = "'0.882"
$ws.Cells.Item(17, 5).Value = "  +2.93%  "
$ws.Cells.Item(18, 4).Value = "'2.246.88"
$ws.Cells.Item(18, 5).Value = "  +2.88%  "
$ws.Cells.Item(19, 4).Value = "'42.741.37"
$ws.Cells.Item(19, 5).Value = "  +4.84%  "
$ws.Cells.Item(20, 4).Value = "'0.0₃0985"
$ws.Cells.Item(20, 5).Value = "  +5.48%  "
$ws.Cells.Item(21, 5).Value = "  +3.52%  "
$ws.Cells.Item(22, 4).Value = "'73.06"
$ws.Cells.Item(22, 5).Value = "  +2.61%  "
$ws.Cells.Item(23, 4).Value = "'235.47"
$ws.Cells.Item(23, 5).Value = "  +2.84%  "
$ws.Cells.Item(24, 4).Value = "'2.05"
$ws.Cells.Item(24, 5).Value = "  +1.07%  "
$ws.Cells.Item(25, 4).Value = "'3.93"
$ws.Cells.Item(25, 5).Value = "  +6.20%  "
$ws.Cells.Item(26, 4).Value = "'11.56"
$ws.Cells.Item(26, 5).Value = "  +2.33%  "
$ws.Cells.Item(27, 4).Value = "'0.998"
$ws.Cells.Item(27, 5).Value = "  -0.18%  "
$ws.Cells.Item(28, 4).Value = "'2.42"
$ws.Cells.Item(28, 5).Value = "  +0.33%  "
$ws.Cells.Item(29, 4).Value = "'3.65"
$ws.Cells.Item(29, 5).Value = "  -1.63%  "
$ws.Cells.Item(30, 4).Value = "'2.21"
$ws.Cells.Item(30, 5).Value = "  +3.78%  "
$ws.Cells.Item(31, 4).Value = "'167.26"
$ws.Cells.Item(31, 5).Value = "  -0.19%  "
$ws.Cells.Item(32, 4).Value = "'20.97"
$ws.Cells.Item(32, 5).Value = "  +4.13%  "
$ws.Cells.Item(33, 4).Value = "'6.46"
$ws.Cells.Item(33, 5).Value = "  +15.61%  "
$ws.Cells.Item(34, 4).Value = "'0.128"
$ws.Cells.Item(34, 5).Value = "  +6.87%  "
$ws.Cells.Item(35, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(35, 4).Value = "'31.38"
$ws.Cells.Item(35, 5).Value = "  +25.34%  "
$ws.Cells.Item(36, 2).Value = "Hedera"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(36, 4).Value = "'0.0790"
$ws.Cells.Item(36, 5).Value = "  +6.92%  "
$ws.Cells.Item(37, 4).Value = "'0.126"
$ws.Cells.Item(37, 5).Value = "  +4.58%  "
$ws.Cells.Item(38, 4).Value = "'4.37"
$ws.Cells.Item(38, 5).Value = "  +7.38%  "
$ws.Cells.Item(39, 4).Value = "'4.71"
$ws.Cells.Item(39, 5).Value = "  +4.10%  "
$ws.Cells.Item(40, 4).Value = "'0.0321"
$ws.Cells.Item(40, 5).Value = "  +8.03%  "
$ws.Cells.Item(41, 4).Value = "'2.31"
$ws.Cells.Item(41, 5).Value = "  +6.89%  "
$ws.Cells.Item(42, 4).Value = "'12.56"
$ws.Cells.Item(42, 5).Value = "  +9.41%  "
$ws.Cells.Item(43, 5).Value = "  +6.01%  "
$ws.Cells.Item(44, 4).Value = "'62.79"
$ws.Cells.Item(44, 5).Value = "  +2.67%  "
$ws.Cells.Item(45, 2).Value = "Algorand"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(45, 4).Value = "'0.201"
$ws.Cells.Item(45, 5).Value = "  +5.05%  "
$ws.Cells.Item(46, 2).Value = "FTXToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(46, 4).Value = "'4.89"
$ws.Cells.Item(46, 5).Value = "  +3.90%  "
$ws.Cells.Item(47, 4).Value = "'8.97"
$ws.Cells.Item(47, 5).Value = "  +5.35%  "
$ws.Cells.Item(48, 5).Value = "  +5.45%  "
$ws.Cells.Item(49, 5).Value = "  +0.20%  "
$ws.Cells.Item(50, 4).Value = "'1.17"
$ws.Cells.Item(50, 5).Value = "  +0.29%  "
$ws.Cells.Item(51, 4).Value = "'1.19"
$ws.Cells.Item(51, 5).Value = "  +4.51%  "

Write-Output "Applied 114 cell updates to Sheet1"
